$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target cluster is now a full 3x3 cross of (ECs, FAPs, MuSCs) sending x receiving clusters (rows 2-10),
# and TPM-derived metrics (cols E-T) are updated throughout.
$rows = @(
    @{ RowNum=2; A="ECs"; B="Nlgn2"; C="Nrxn1"; D="ECs"; E=3; F=1; G=1.164555333333333; H=3.493666; I=0.05923394707027321; J=0.05923394707027322; K=1; L=0.3333333333333333; M=0.002858666666666667; N=0.008576; O=0.002669495535069502; P=0.002669495535069501; Q=0.003329075512888889; R=0.029961679616; S=0.0001581247572286375; T=0.0001581247572286375 },
    @{ RowNum=3; A="ECs"; B="Nlgn2"; C="Nrxn1"; D="FAPs"; E=3; F=1; G=1.164555333333333; H=3.493666; I=0.05923394707027321; J=0.05923394707027322; K=2; L=0.6666666666666666; M=0.07823633333333334; N=0.234709; O=0.07305907503971872; P=0.0730590750397187; Q=0.09111053924377779; R=0.819994853194; S=0.004327577383905817; T=0.004327577383905817 },
    @{ RowNum=4; A="ECs"; B="Nlgn2"; C="Nrxn1"; D="MuSCs"; E=3; F=1; G=1.164555333333333; H=3.493666; I=0.05923394707027321; J=0.05923394707027322; K=3; L=1; M=0.989769; N=2.969307; O=0.9242714294252118; P=0.9242714294252118; Q=1.152640767718; R=10.373766909462; S=0.05474824492913876; T=0.05474824492913877 },
    @{ RowNum=5; A="FAPs"; B="Nlgn2"; C="Nrxn1"; D="ECs"; E=3; F=1; G=10.79119133333334; H=32.373574; I=0.5488831985632208; J=0.5488831985632209; K=1; L=0.3333333333333333; M=0.002858666666666667; N=0.008576; O=0.002669495535069502; P=0.002669495535069501; Q=0.03084841895822223; R=0.2776357706240001; S=0.001465241247839185; T=0.001465241247839185 },
    @{ RowNum=6; A="FAPs"; B="Nlgn2"; C="Nrxn1"; D="FAPs"; E=3; F=1; G=10.79119133333334; H=32.373574; I=0.5488831985632208; J=0.5488831985632209; K=2; L=0.6666666666666666; M=0.07823633333333334; N=0.234709; O=0.07305907503971872; P=0.0730590750397187; Q=0.8442632422184446; R=7.598369179966001; S=0.04010089879187118; T=0.04010089879187118 },
    @{ RowNum=7; A="FAPs"; B="Nlgn2"; C="Nrxn1"; D="MuSCs"; E=3; F=1; G=10.79119133333334; H=32.373574; I=0.5488831985632208; J=0.5488831985632209; K=3; L=1; M=0.989769; N=2.969307; O=0.9242714294252118; P=0.9242714294252118; Q=10.680786654802; R=96.12707989321802; S=0.5073170585235104; T=0.5073170585235105 },
    @{ RowNum=8; A="MuSCs"; B="Nlgn2"; C="Nrxn1"; D="ECs"; E=3; F=1; G=7.704522333333333; H=23.113567; I=0.391882854366506; J=0.3918828543665061; K=1; L=0.3333333333333333; M=0.002858666666666667; N=0.008576; O=0.002669495535069502; P=0.002669495535069501; Q=0.02202466117688889; R=0.198221950592; S=0.00104612953000168; T=0.00104612953000168 },
    @{ RowNum=9; A="MuSCs"; B="Nlgn2"; C="Nrxn1"; D="FAPs"; E=3; F=1; G=7.704522333333333; H=23.113567; I=0.391882854366506; J=0.3918828543665061; K=2; L=0.6666666666666666; M=0.07823633333333334; N=0.234709; O=0.07305907503971872; P=0.0730590750397187; Q=0.6027735774447778; R=5.424962197003; S=0.02863059886394172; T=0.02863059886394172 },
    @{ RowNum=10; A="MuSCs"; B="Nlgn2"; C="Nrxn1"; D="MuSCs"; E=3; F=1; G=7.704522333333333; H=23.113567; I=0.391882854366506; J=0.3918828543665061; K=3; L=1; M=0.989769; N=2.969307; O=0.9242714294252118; P=0.9242714294252118; Q=7.625697365341; R=68.631276288069; S=0.3622061259725626; T=0.3622061259725626 }
)

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")
foreach ($row in $rows) {
    foreach ($col in $cols) {
        $ws.Range("$col$($row.RowNum)").Value = $row[$col]
    }
}